$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph near the top of the document
#    (right after the "Play Bounty Showdown Free ..." Heading1 paragraph).
#    It currently reads: "Meta description" (bold run) + ": Read our review
#    of ..." (plain run). The whole paragraph, including its paragraph mark,
#    is deleted.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new paragraph just before the final paragraph of the document,
#    containing the bold heading text "Play Bounty Showdown Free - A Wild
#    West Video Slot Game". We splice in raw WordOpenXML (matching the
#    leading empty-run style used throughout this document) at the very
#    start of the final paragraph, then split it into its own paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$startPos = $lastPara.Range.Start
$newHeadingText = "Play Bounty Showdown Free - A Wild West Video Slot Game"

$insertPoint = $d.Range($startPos, $startPos)
$openXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $newHeadingText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($openXml)

# Split the document right after the text we just inserted, so it becomes
# its own paragraph instead of sharing one with the original last paragraph.
$breakPos = $startPos + $newHeadingText.Length
$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3. Replace the text of the (new) final paragraph -- previously the
#    feature-image prompt -- with the meta-description text, keeping its
#    italic formatting intact.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$oldPrompt = "Create an eye-catching feature image for Bounty Showdown that incorporates the Wild West theme and the Mayan element. The cartoon-style image should feature a happy Maya warrior wearing glasses. The warrior could be sitting on top of a horse, holding a lasso, or standing in front of a wooden saloon. The background of the image should showcase the Wild West landscape, including red rock formations, a Western town, or a dusty desert. The text " + [char]34 + "Bounty Showdown" + [char]34 + " should be prominently displayed in an Old West-style font. The overall feel of the image should be fun and vibrant, with bright colors and playful details."
$newPrompt = "Read our review of Bounty Showdown, a Wild West themed video slot game with bonus features. Play it for free on desktop, tablet, and mobile devices."
$found = $find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)
Write-Output ("replaced feature-image prompt: " + $found)
